$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update runs/balls/fours/sixes (columns C,D,E,F) for rows 2,4,5,7,9
# to reflect the corrected match-by-match activity log.

$ws.Range("C2").Value = "12"
$ws.Range("D2").Value = "7"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "1"

$ws.Range("C4").Value = "15"
$ws.Range("D4").Value = "15"
$ws.Range("E4").Value = "2"
$ws.Range("F4").Value = "0"

$ws.Range("C5").Value = "9"
$ws.Range("D5").Value = "10"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "0"

$ws.Range("C7").Value = "15"
$ws.Range("D7").Value = "7"
$ws.Range("E7").Value = "1"
$ws.Range("F7").Value = "1"

$ws.Range("C9").Value = "2"
$ws.Range("D9").Value = "3"
$ws.Range("E9").Value = "0"
$ws.Range("F9").Value = "0"
